# ---------------------------------------------------------------------------
# Fix IXL APZ computation template (R_ZSM_3 verification sheet)
#  - rename "IXL APZ Limit Downstream" -> "Downstream Limit"
#  - rename "IXL APZ Limit Upstream"   -> "Upstream Limit"
#  - add "Type" and "Direction" columns right after "Signal Name"
#  - "Related CBTC Direction Zone" header changes colour (orange -> pink)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R_ZSM_3")

# Colours used across the header row (OLE/BGR long values == RGB(r,g,b))
$colOrange = 10079487   # FFCC99 - Type / Direction (new)
$colPink   = 13408767   # FF99CC - Related CBTC Direction Zone (new colour)

# 1) Rename the merged "IXL APZ Limit Downstream" header (D1:G1) to
#    "Downstream Limit", and "IXL APZ Limit Upstream" (H1:K1) to
#    "Upstream Limit" -- done first (this is how the shared strings end up
#    ordered in the saved file).
$ws.Range("D1").Value = "Downstream Limit"
$ws.Range("H1").Value = "Upstream Limit"

# 2) Insert two new blank columns before the old column B so that the old
#    B..O columns (Related CBTC Direction Zone .. Comments) shift to D..Q.
$ws.Columns("B:C").Insert()

# 3) New "Type" / "Direction" headers (B1:B2 and C1:C2), styled like the old
#    "Related CBTC Direction Zone" header (orange, bold Arial 10, centred,
#    wrapped, boxed).
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Direction"

foreach ($hdr in @($ws.Range("B1:B2"), $ws.Range("C1:C2"))) {
    $hdr.Interior.Pattern = -4124   # xlSolid
    $hdr.Interior.Color = $colOrange
    $hdr.Font.Name = "Arial"
    $hdr.Font.Size = 10
    $hdr.Font.Bold = $true
    $hdr.HorizontalAlignment = -4108  # xlCenter
    $hdr.VerticalAlignment = -4108    # xlCenter
    $hdr.WrapText = $true
    $hdr.NumberFormat = "General"
    $hdr.Borders.LineStyle = 1
    $hdr.Borders.Weight = 2
    $hdr.Merge() | Out-Null
}

# 4) "Related CBTC Direction Zone" header (now D1:D2) switches from orange to pink.
$relHdr = $ws.Range("D1:D2")
$relHdr.Interior.Pattern = -4124
$relHdr.Interior.Color = $colPink

# 5) Fix the AutoFilter range (now spans A2:Q2 instead of A2:O2).
$ws.AutoFilterMode = $false
$ws.Range("A2:Q2").AutoFilter() | Out-Null

# 6) Fix the conditional-formatting ranges that used to point at N/O so that
#    they follow the columns to their new P/Q positions (column insertion
#    does not retarget conditional formatting automatically).
$fcs = $ws.Cells.FormatConditions
for ($i = 1; $i -le 3; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("P1:P1048576")) | Out-Null
}
for ($i = 4; $i -le 6; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("P1:Q1")) | Out-Null
}

Write-Host "Done."
